$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Ana Milic"
$ws.Range("B8").Value = "Kombinovan frižider"
$ws.Range("C8").Value = "Samsung"
$ws.Range("D8").Value = "RF-850"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "22222"
$ws.Range("F8").Value = "ne hladi zamrzivač"
